# Daily attendance processing - swap the order of the first two
# comma-separated entries in the "Recorded By" column (column G) so that
# the human/email entry is listed before the "System" entry.
#
# Rule (derived from the target diff):
#   - Cells with a single value are left untouched.
#   - Cells whose value is exactly "admin@admin.com, System" are left
#     untouched.
#   - All other multi-value cells have their first two comma-separated
#     tokens swapped; any additional trailing tokens (e.g. a lowercase
#     "system" duplicate) keep their position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)   # Column G = "Recorded By"
    $value = $cell.Value2

    if ($null -eq $value) { continue }
    if (-not ($value -is [string])) { continue }
    if ($value -eq "") { continue }
    if ($value -eq "admin@admin.com, System") { continue }

    $parts = $value -split ", "
    if ($parts.Count -lt 2) { continue }

    $first = $parts[0]
    $second = $parts[1]
    $parts[0] = $second
    $parts[1] = $first

    $cell.Value2 = [string]::Join(", ", $parts)
}
